$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Software Engineering Placement 2025 - (Cloud DevOps)"
$ws.Range("B14").Value = "Hewlett Packard Enterprise"
$ws.Range("C14").Value = "Software development methodologies, principles, practices, software development lifecycle, C, C#, C++, Java, Python, Go, Computer architecture, Concurrent programming/multi-threading, Embedded systems, Linux, Networks and communications, Cloud Microservice development, Operating Systems"
$ws.Range("D14").Value = "Bristol, UK"

$ws.Range("A15").Value = "2026 Technology (GOTO) Industrial Placement"
$ws.Range("B15").Value = "UBS"
$ws.Range("C15").Value = "Financial Management, Information Technology, Computer Science, Software Engineering, Investment Banking"
$ws.Range("D15").Value = "London"
